$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

# New rows to append: Question, Tag_1, Tag_2, Tag_3, Tag_4
$data = @(
    @("ทำไมสินค้ายี่ห้อ ""Apple"" ถึงราคาแพง", "Apple", "สินค้า", "ของแพง", $null),
    @("ท้องเสียมา 2 วันติดแล้ว กินยาตัวไหนดี", "ท้องเสีย", "เจ็บป่วย", "ยา", $null),
    @("เครื่องซักผ้าเสีย ซักผ้าอยู่ดี ๆ ก็ดับไปเลย", "เครื่องซักผ้า", "เครื่องใช้ไฟฟ้า", "ซักผ้า", "งานบ้าน"),
    @("แม่ผัวน่ารำคาญมาก หาเรื่องบ่นเราได้ตลอดเลย ทำยังไงได้บ้างคะ", "แม่ผัว", "ครอบครัว", "ปัญหาครอบครัว", $null),
    @("กลิ่นเต่าแฟนเหม็นมาก รับไม่ได้เลย สามารถทำยังไงได้บ้าง", "กลิ่นเต่า", "ร่างกาย", "กลิ่นรักแร้", $null)
)

$startRow = 10
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        if ($null -ne $rowValues[$j]) {
            $ws.Cells.Item($rowIndex, $j + 1).Value = $rowValues[$j]
        }
    }
}
